# Added 2050 FBP PlusCrossing_22 run
# Insert a new data row at row 137 of the all_runs sheet (pushing the
# existing rows 137-147 down to 138-148) and populate it with the new
# 2050 FBP PlusCrossing run, matching the style/formatting that Excel
# carries down automatically from the row above (135/136).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# Insert a new row above the current row 137 (the first EIR Alt1 row),
# shifting all the rows below it down by one.
$ws.Rows(137).Insert()

# Populate the newly-inserted row with the new run's data.
$ws.Range("A137").Value() = "RTP2021"
$ws.Range("B137").Value() = 2050
$ws.Range("C137").Value() = "2050_TM152_FBP_PlusCrossing_22"
$ws.Range("D137").Value() = "FinalBlueprint"
$ws.Range("E137").Value() = "Plus"
$ws.Range("F137").Value() = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION"'
$ws.Range("G137").Value() = "run182"
$ws.Range("H137").Value() = "current"

# Match the author's final cursor position / selection on the sheet.
$ws.Range("C137").Select() | Out-Null

